# "Generate Report for Archive"
# The localization status report is regenerated: the pending "Ready for
# handoff" status becomes "In Translation" everywhere it is used (the
# per-language Overview columns and each language sheet's Status column),
# and the Status column is narrower to match the shorter text after the
# report re-auto-fit its columns.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn/de-de status columns (E, F) for both rows ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# Narrow columns E and F to the new auto-fit width
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

# --- Per-language sheets: Status column (C) for both data rows ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C1").ColumnWidth = 12.5
